# Update Clec11a-Itga11 LR-pairs sheet with new TPM-derived values.
# The sending/target cluster set grows from {ECs, FAPs, MuSCs} to
# {ECs, FAPs, MuSCs, Resolving-Mac}, producing a full 4x4 cross (16 rows)
# instead of the previous 4x3 (12 rows), and all numeric columns are
# refreshed with the new TPM-based statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Clec11a"
$ws.Cells.Item(2,3).Value = "Itga11"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.152918
$ws.Cells.Item(2,8).Value = 0.458754
$ws.Cells.Item(2,9).Value = 0.007158013163202275
$ws.Cells.Item(2,10).Value = 0.007158013163202275
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.152959
$ws.Cells.Item(2,14).Value = 0.458877
$ws.Cells.Item(2,15).Value = 0.004761500378002596
$ws.Cells.Item(2,16).Value = 0.004761500378002596
$ws.Cells.Item(2,17).Value = 0.023390184362
$ws.Cells.Item(2,18).Value = 0.210511659258
$ws.Cells.Item(2,19).Value = 0.00003408288238233519
$ws.Cells.Item(2,20).Value = 0.00003408288238233519
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Clec11a"
$ws.Cells.Item(3,3).Value = "Itga11"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.152918
$ws.Cells.Item(3,8).Value = 0.458754
$ws.Cells.Item(3,9).Value = 0.007158013163202275
$ws.Cells.Item(3,10).Value = 0.007158013163202275
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 31.603318
$ws.Cells.Item(3,14).Value = 94.809954
$ws.Cells.Item(3,15).Value = 0.9837878817404418
$ws.Cells.Item(3,16).Value = 0.9837878817404418
$ws.Cells.Item(3,17).Value = 4.832716181924
$ws.Cells.Item(3,18).Value = 43.49444563731601
$ws.Cells.Item(3,19).Value = 0.007041966607296966
$ws.Cells.Item(3,20).Value = 0.007041966607296966
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Clec11a"
$ws.Cells.Item(4,3).Value = "Itga11"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.152918
$ws.Cells.Item(4,8).Value = 0.458754
$ws.Cells.Item(4,9).Value = 0.007158013163202275
$ws.Cells.Item(4,10).Value = 0.007158013163202275
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.2671263333333334
$ws.Cells.Item(4,14).Value = 0.8013790000000001
$ws.Cells.Item(4,15).Value = 0.008315444904458803
$ws.Cells.Item(4,16).Value = 0.008315444904458805
$ws.Cells.Item(4,17).Value = 0.04084842464066667
$ws.Cells.Item(4,18).Value = 0.367635821766
$ws.Cells.Item(4,19).Value = 0.0000595220640839994
$ws.Cells.Item(4,20).Value = 0.00005952206408399941
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Clec11a"
$ws.Cells.Item(5,3).Value = "Itga11"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.152918
$ws.Cells.Item(5,8).Value = 0.458754
$ws.Cells.Item(5,9).Value = 0.007158013163202275
$ws.Cells.Item(5,10).Value = 0.007158013163202275
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.1007146666666667
$ws.Cells.Item(5,14).Value = 0.302144
$ws.Cells.Item(5,15).Value = 0.00313517297709673
$ws.Cells.Item(5,16).Value = 0.00313517297709673
$ws.Cells.Item(5,17).Value = 0.01540108539733333
$ws.Cells.Item(5,18).Value = 0.138609768576
$ws.Cells.Item(5,19).Value = 0.00002244160943897446
$ws.Cells.Item(5,20).Value = 0.00002244160943897446
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Clec11a"
$ws.Cells.Item(6,3).Value = "Itga11"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 20.82581466666667
$ws.Cells.Item(6,8).Value = 62.47744400000001
$ws.Cells.Item(6,9).Value = 0.9748457050079848
$ws.Cells.Item(6,10).Value = 0.9748457050079848
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.152959
$ws.Cells.Item(6,14).Value = 0.458877
$ws.Cells.Item(6,15).Value = 0.004761500378002596
$ws.Cells.Item(6,16).Value = 0.004761500378002596
$ws.Cells.Item(6,17).Value = 3.185495785598667
$ws.Cells.Item(6,18).Value = 28.669462070388
$ws.Cells.Item(6,19).Value = 0.004641728192889727
$ws.Cells.Item(6,20).Value = 0.004641728192889727
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Clec11a"
$ws.Cells.Item(7,3).Value = "Itga11"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 20.82581466666667
$ws.Cells.Item(7,8).Value = 62.47744400000001
$ws.Cells.Item(7,9).Value = 0.9748457050079848
$ws.Cells.Item(7,10).Value = 0.9748457050079848
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 31.603318
$ws.Cells.Item(7,14).Value = 94.809954
$ws.Cells.Item(7,15).Value = 0.9837878817404418
$ws.Cells.Item(7,16).Value = 0.9837878817404418
$ws.Cells.Item(7,17).Value = 658.1648435197308
$ws.Cells.Item(7,18).Value = 5923.483591677576
$ws.Cells.Item(7,19).Value = 0.959041391153573
$ws.Cells.Item(7,20).Value = 0.959041391153573
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Clec11a"
$ws.Cells.Item(8,3).Value = "Itga11"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 20.82581466666667
$ws.Cells.Item(8,8).Value = 62.47744400000001
$ws.Cells.Item(8,9).Value = 0.9748457050079848
$ws.Cells.Item(8,10).Value = 0.9748457050079848
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.2671263333333334
$ws.Cells.Item(8,14).Value = 0.8013790000000001
$ws.Cells.Item(8,15).Value = 0.008315444904458803
$ws.Cells.Item(8,16).Value = 0.008315444904458805
$ws.Cells.Item(8,17).Value = 5.563123510586223
$ws.Cells.Item(8,18).Value = 50.06811159527601
$ws.Cells.Item(8,19).Value = 0.008106275750342198
$ws.Cells.Item(8,20).Value = 0.0081062757503422
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Clec11a"
$ws.Cells.Item(9,3).Value = "Itga11"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 20.82581466666667
$ws.Cells.Item(9,8).Value = 62.47744400000001
$ws.Cells.Item(9,9).Value = 0.9748457050079848
$ws.Cells.Item(9,10).Value = 0.9748457050079848
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.1007146666666667
$ws.Cells.Item(9,14).Value = 0.302144
$ws.Cells.Item(9,15).Value = 0.00313517297709673
$ws.Cells.Item(9,16).Value = 0.00313517297709673
$ws.Cells.Item(9,17).Value = 2.097464982215111
$ws.Cells.Item(9,18).Value = 18.877184839936
$ws.Cells.Item(9,19).Value = 0.003056309911179844
$ws.Cells.Item(9,20).Value = 0.003056309911179845
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Clec11a"
$ws.Cells.Item(10,3).Value = "Itga11"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.200477
$ws.Cells.Item(10,8).Value = 0.601431
$ws.Cells.Item(10,9).Value = 0.009384225564807953
$ws.Cells.Item(10,10).Value = 0.009384225564807953
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.152959
$ws.Cells.Item(10,14).Value = 0.458877
$ws.Cells.Item(10,15).Value = 0.004761500378002596
$ws.Cells.Item(10,16).Value = 0.004761500378002596
$ws.Cells.Item(10,17).Value = 0.030664761443
$ws.Cells.Item(10,18).Value = 0.275982852987
$ws.Cells.Item(10,19).Value = 0.0000446829935740947
$ws.Cells.Item(10,20).Value = 0.0000446829935740947
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Clec11a"
$ws.Cells.Item(11,3).Value = "Itga11"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.200477
$ws.Cells.Item(11,8).Value = 0.601431
$ws.Cells.Item(11,9).Value = 0.009384225564807953
$ws.Cells.Item(11,10).Value = 0.009384225564807953
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 31.603318
$ws.Cells.Item(11,14).Value = 94.809954
$ws.Cells.Item(11,15).Value = 0.9837878817404418
$ws.Cells.Item(11,16).Value = 0.9837878817404418
$ws.Cells.Item(11,17).Value = 6.335738382686001
$ws.Cells.Item(11,18).Value = 57.02164544417401
$ws.Cells.Item(11,19).Value = 0.009232087390176918
$ws.Cells.Item(11,20).Value = 0.009232087390176918
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Clec11a"
$ws.Cells.Item(12,3).Value = "Itga11"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.200477
$ws.Cells.Item(12,8).Value = 0.601431
$ws.Cells.Item(12,9).Value = 0.009384225564807953
$ws.Cells.Item(12,10).Value = 0.009384225564807953
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.2671263333333334
$ws.Cells.Item(12,14).Value = 0.8013790000000001
$ws.Cells.Item(12,15).Value = 0.008315444904458803
$ws.Cells.Item(12,16).Value = 0.008315444904458805
$ws.Cells.Item(12,17).Value = 0.05355268592766668
$ws.Cells.Item(12,18).Value = 0.4819741733490001
$ws.Cells.Item(12,19).Value = 0.00007803401065517433
$ws.Cells.Item(12,20).Value = 0.00007803401065517434
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Clec11a"
$ws.Cells.Item(13,3).Value = "Itga11"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.200477
$ws.Cells.Item(13,8).Value = 0.601431
$ws.Cells.Item(13,9).Value = 0.009384225564807953
$ws.Cells.Item(13,10).Value = 0.009384225564807953
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.1007146666666667
$ws.Cells.Item(13,14).Value = 0.302144
$ws.Cells.Item(13,15).Value = 0.00313517297709673
$ws.Cells.Item(13,16).Value = 0.00313517297709673
$ws.Cells.Item(13,17).Value = 0.02019097422933333
$ws.Cells.Item(13,18).Value = 0.181718768064
$ws.Cells.Item(13,19).Value = 0.00002942117040176619
$ws.Cells.Item(13,20).Value = 0.0000294211704017662
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Clec11a"
$ws.Cells.Item(14,3).Value = "Itga11"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 0.6666666666666666
$ws.Cells.Item(14,7).Value = 0.183981
$ws.Cells.Item(14,8).Value = 0.5519430000000001
$ws.Cells.Item(14,9).Value = 0.008612056264005009
$ws.Cells.Item(14,10).Value = 0.008612056264005009
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.152959
$ws.Cells.Item(14,14).Value = 0.458877
$ws.Cells.Item(14,15).Value = 0.004761500378002596
$ws.Cells.Item(14,16).Value = 0.004761500378002596
$ws.Cells.Item(14,17).Value = 0.02814154977900001
$ws.Cells.Item(14,18).Value = 0.2532739480110001
$ws.Cells.Item(14,19).Value = 0.00004100630915643948
$ws.Cells.Item(14,20).Value = 0.00004100630915643948
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Clec11a"
$ws.Cells.Item(15,3).Value = "Itga11"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 0.6666666666666666
$ws.Cells.Item(15,7).Value = 0.183981
$ws.Cells.Item(15,8).Value = 0.5519430000000001
$ws.Cells.Item(15,9).Value = 0.008612056264005009
$ws.Cells.Item(15,10).Value = 0.008612056264005009
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 31.603318
$ws.Cells.Item(15,14).Value = 94.809954
$ws.Cells.Item(15,15).Value = 0.9837878817404418
$ws.Cells.Item(15,16).Value = 0.9837878817404418
$ws.Cells.Item(15,17).Value = 5.814410048958002
$ws.Cells.Item(15,18).Value = 52.32969044062201
$ws.Cells.Item(15,19).Value = 0.008472436589394991
$ws.Cells.Item(15,20).Value = 0.008472436589394991
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Clec11a"
$ws.Cells.Item(16,3).Value = "Itga11"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 0.6666666666666666
$ws.Cells.Item(16,7).Value = 0.183981
$ws.Cells.Item(16,8).Value = 0.5519430000000001
$ws.Cells.Item(16,9).Value = 0.008612056264005009
$ws.Cells.Item(16,10).Value = 0.008612056264005009
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.2671263333333334
$ws.Cells.Item(16,14).Value = 0.8013790000000001
$ws.Cells.Item(16,15).Value = 0.008315444904458803
$ws.Cells.Item(16,16).Value = 0.008315444904458805
$ws.Cells.Item(16,17).Value = 0.04914616993300001
$ws.Cells.Item(16,18).Value = 0.4423155293970001
$ws.Cells.Item(16,19).Value = 0.00007161307937743296
$ws.Cells.Item(16,20).Value = 0.00007161307937743298
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Clec11a"
$ws.Cells.Item(17,3).Value = "Itga11"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 0.6666666666666666
$ws.Cells.Item(17,7).Value = 0.183981
$ws.Cells.Item(17,8).Value = 0.5519430000000001
$ws.Cells.Item(17,9).Value = 0.008612056264005009
$ws.Cells.Item(17,10).Value = 0.008612056264005009
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.1007146666666667
$ws.Cells.Item(17,14).Value = 0.302144
$ws.Cells.Item(17,15).Value = 0.00313517297709673
$ws.Cells.Item(17,16).Value = 0.00313517297709673
$ws.Cells.Item(17,17).Value = 0.018529585088
$ws.Cells.Item(17,18).Value = 0.166766265792
$ws.Cells.Item(17,19).Value = 0.00002700028607614512
$ws.Cells.Item(17,20).Value = 0.00002700028607614513
